$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: stash the row that is about to be pushed out of the 2..23 block ---
# Row 23 ("Altona Pier, ...") currently has a blank LGA_NAME20 (column D). Once the
# rows below row 2 shift down by one, this row's data is relocated to the very end
# of the sheet (new last row) with its LGA_NAME20 filled in.
$tailA = $ws.Range("A23").Value2
$tailB = $ws.Range("B23").Value2
$tailC = $ws.Range("C23").Value2

# --- Step 2: shift rows 2..22 down into rows 3..23 (iterate bottom-up so we never
# clobber a source row before it has been read) ---
for ($r = 22; $r -ge 2; $r--) {
    $ws.Cells.Item($r + 1, 1).Value = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r + 1, 2).Value = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r + 1, 3).Value = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($r + 1, 4).Value = $ws.Cells.Item($r, 4).Value2
}

# --- Step 3: write the brand-new first data row ---
$ws.Range("A2").Value = "Shelton Court, Noble Park North VIC 3170, Australia"
$ws.Range("B2").Value = -37.9408701
$ws.Range("C2").Value = 145.1962728
$ws.Range("D2").Value = "Greater Dandenong (C)"

# --- Step 4: append the stashed row at the bottom of the sheet, with its
# previously-missing LGA_NAME20 now populated ---
$lastRow = $ws.UsedRange.Rows.Count + 1
$ws.Cells.Item($lastRow, 1).Value = $tailA
$ws.Cells.Item($lastRow, 2).Value = $tailB
$ws.Cells.Item($lastRow, 3).Value = $tailC
$ws.Cells.Item($lastRow, 4).Value = "Hobsons Bay (C)"
